$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 39900
$ws.Range("J57").Value = 40000
$ws.Range("L57").Value = 120000
$ws.Range("N57").Value = -120998

$ws.Range("H58").Value = 38462.184
$ws.Range("I58").Value = 392.23077
$ws.Range("J58").Value = 73812.86
$ws.Range("K58").Value = 1176.69231
$ws.Range("L58").Value = 221438.58
$ws.Range("M58").Value = -1026.69231
$ws.Range("N58").Value = -221738.58

$ws.Range("H76").Value = 2703.8462
$ws.Range("I76").Value = 2703.8462
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2703.8462
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2388.8462
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 2703.8462
$ws.Range("I79").Value = 2703.8462
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2703.8462
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1611.8462

$ws.Range("H116").Value = 4020.3572
$ws.Range("I116").Value = 3634.1428
$ws.Range("K116").Value = 3634.1428
$ws.Range("M116").Value = -192.1428000000001

$ws.Range("H132").Value = 6252752.5
$ws.Range("I132").Value = 8002738
$ws.Range("J132").Value = 2805.7144
$ws.Range("K132").Value = 24008214
$ws.Range("L132").Value = 8417.143199999999
$ws.Range("M132").Value = -24005684
$ws.Range("N132").Value = -13477.1432

$ws.Range("H137").Value = 7400
$ws.Range("I137").Value = 10061.2
$ws.Range("J137").Value = 4073.5
$ws.Range("K137").Value = 30183.6
$ws.Range("L137").Value = 12220.5
$ws.Range("M137").Value = -27633.6
$ws.Range("N137").Value = -17320.5

$ws.Range("H138").Value = 1751.0405
$ws.Range("I138").Value = 1085.3889
$ws.Range("K138").Value = 3256.1667
$ws.Range("M138").Value = 1883.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7144269.5
$ws.Range("I2").Value = 13159369
$ws.Range("J2").Value = 1338.875
$ws.Range("K2").Value = 13159369
$ws.Range("L2").Value = 1338.875
$ws.Range("M2").Value = -13159256
$ws.Range("N2").Value = -1564.875

$ws.Range("H63").Value = 2634.1177
$ws.Range("I63").Value = 1999.1666
$ws.Range("J63").Value = 4158
$ws.Range("K63").Value = 1999.1666
$ws.Range("L63").Value = 4158
$ws.Range("M63").Value = -1313.1666
$ws.Range("N63").Value = -5530

$ws.Range("H66").Value = 2634.1177
$ws.Range("I66").Value = 1999.1666
$ws.Range("J66").Value = 4158
$ws.Range("K66").Value = 9995.833000000001
$ws.Range("L66").Value = 20790
$ws.Range("M66").Value = -6563.833000000001
$ws.Range("N66").Value = -27654

$ws.Range("H102").Value = 2960
$ws.Range("I102").Value = 2954.2856
$ws.Range("K102").Value = 2954.2856
$ws.Range("M102").Value = -1332.2856

$ws.Range("H116").Value = 7144269.5
$ws.Range("I116").Value = 13159369
$ws.Range("J116").Value = 1338.875
$ws.Range("K116").Value = 13159369
$ws.Range("L116").Value = 1338.875
$ws.Range("M116").Value = -13157075
$ws.Range("N116").Value = -5926.875

$ws.Range("H123").Value = 29000
$ws.Range("J123").Value = 29000
$ws.Range("L123").Value = 29000
$ws.Range("N123").Value = -38800

$ws.Range("H132").Value = 2179.9216
$ws.Range("I132").Value = 2021.25
$ws.Range("K132").Value = 6063.75
$ws.Range("M132").Value = -3533.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7144269.5
$ws.Range("I3").Value = 13159369
$ws.Range("J3").Value = 1338.875
$ws.Range("K3").Value = 13159369
$ws.Range("L3").Value = 1338.875
$ws.Range("M3").Value = -13159255
$ws.Range("N3").Value = -1566.875

$ws.Range("H105").Value = 1283.5264
$ws.Range("I105").Value = 1271.2693
$ws.Range("J105").Value = 1310.0834
$ws.Range("K105").Value = 1271.2693
$ws.Range("L105").Value = 1310.0834
$ws.Range("M105").Value = 475.7307000000001
$ws.Range("N105").Value = -4804.0834

$ws.Range("H134").Value = 1928.3
$ws.Range("I134").Value = 1503.4736
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 4510.4208
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -1975.4208
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3790.4707
$ws.Range("I31").Value = 2612.8696
$ws.Range("K31").Value = 2612.8696
$ws.Range("M31").Value = -2317.8696

$ws.Range("H34").Value = 3790.4707
$ws.Range("I34").Value = 2612.8696
$ws.Range("K34").Value = 2612.8696
$ws.Range("M34").Value = -2410.8696

$ws.Range("H132").Value = 2659.6365
$ws.Range("I132").Value = 2739.3845
$ws.Range("J132").Value = 2544.4443
$ws.Range("K132").Value = 8218.1535
$ws.Range("L132").Value = 7633.3329
$ws.Range("M132").Value = -5688.1535
$ws.Range("N132").Value = -12693.3329

$ws.Range("H134").Value = 6699.5
$ws.Range("I134").Value = 7399.385
$ws.Range("J134").Value = 3666.6667
$ws.Range("K134").Value = 22198.155
$ws.Range("L134").Value = 11000.0001
$ws.Range("M134").Value = -19663.155
$ws.Range("N134").Value = -16070.0001

$ws.Range("H139").Value = 40000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 40000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 40000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2087.7778
$ws.Range("I116").Value = 1216.6666
$ws.Range("J116").Value = 3830
$ws.Range("K116").Value = 3649.9998
$ws.Range("L116").Value = 11490
$ws.Range("M116").Value = -207.9998000000001
$ws.Range("N116").Value = -18374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H136").Value = 2502.0645
$ws.Range("I136").Value = 2114
$ws.Range("J136").Value = 4520
$ws.Range("K136").Value = 6342
$ws.Range("L136").Value = 13560
$ws.Range("M136").Value = -3792
$ws.Range("N136").Value = -18660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 23762.25
$ws.Range("J39").Value = 23762.25
$ws.Range("L39").Value = 23762.25
$ws.Range("N39").Value = -24588.25
